# ResumenResultadosaValidarMonitoreo.xlsx -- add the "no cumple fecha de
# entrega" authorization column, reword a couple of the existing headers,
# and refresh the header-row look (bold) plus a few column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates -------------------------------------------------
# R1 used to read "Se corre reglas de validación"; that question now moves
# all the way to the new last column (U1), and R1 takes over the (slightly
# re-spaced) "Cumple con todos los criterios..." wording that used to live
# in T1.
$ws.Range("R1").Value = "Cumple con todos los criterios para aplicar reglas?  (SI/NO)"

# S1: drop the trailing period, broaden "(SI)" to "(SI/NO)".
$ws.Range("S1").Value = "Autorización de reglas cuando esté incompleto (SI/NO)"

# T1: brand new header (new business rule around delivery-date compliance).
$ws.Range("T1").Value = "Autorización de reglas cuando no cumple fecha de entrega (SI/NO)"

# --- New column U ----------------------------------------------------------
$ws.Range("U1").Value = "Se corre reglas de validación"

# --- Header row formatting ------------------------------------------------
# Center align + wrap the whole header row (A1:U1); keep the original
# template fill, just make the font bold.
$headerRange = $ws.Range("A1:U1")
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4108    # xlCenter
$headerRange.WrapText = $true
$headerRange.Font.Bold = $true

# The brand-new U1 cell has no inherited fill yet -- copy the fully-resolved
# header format (bold font + fill + alignment) from A1 so it matches the
# rest of the header row exactly.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("U1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Column P ("Obs. Condicionantes") keeps a centered + wrapped look instead
# of the old left/top alignment.
$ws.Columns.Item(16).HorizontalAlignment = -4108
$ws.Columns.Item(16).VerticalAlignment = -4108
$ws.Columns.Item(16).WrapText = $true

# --- Column widths (new/changed columns) ----------------------------------
$ws.Columns.Item(15).ColumnWidth = 13.666666666666666   # O  ~14.44
$ws.Columns.Item(16).ColumnWidth = 95.16666666666667    # P  96
$ws.Columns.Item(18).ColumnWidth = 23.0                 # R  ~23.89
$ws.Columns.Item(19).ColumnWidth = 20.666666666666668   # S  ~21.44
$ws.Columns.Item(20).ColumnWidth = 26.333333333333332   # T  ~27.11

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1  # xlPortrait

# --- View state -------------------------------------------------------------
$ws.Range("D6:E6").Select()

$wb.Save()
